$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.492.91'
$ws.Range("E2").Value = '  -1.02%  '
$ws.Range("D3").Value = '2.063.93'
$ws.Range("E3").Value = '  -0.99%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.67'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.623'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.53'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.96%  '
$ws.Range("E9").Value = '  -2.09%  '
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("E11").Value = '  +1.52%  '
$ws.Range("D12").Value = '2.366.86'
$ws.Range("E12").Value = '  -1.03%  '
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("E15").Value = '  -1.72%  '
$ws.Range("E16").Value = '  +0.13%  '
$ws.Range("D17").Value = '2.064.97'
$ws.Range("E17").Value = '  -0.75%  '
$ws.Range("D18").Value = '37.487.69'
$ws.Range("E18").Value = '  -0.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.80'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.39'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("E25").Value = '  -2.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.26'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.00%  '
$ws.Range("E28").Value = '  -4.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("E30").Value = '  -4.23%  '
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("E32").Value = '  -3.59%  '
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("E34").Value = '  -2.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.87%  '
$ws.Range("E36").Value = '  +0.36%  '
$ws.Range("E37").Value = '  -3.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("E39").Value = '  -1.73%  '
$ws.Range("E40").Value = '  +3.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("E42").Value = '  +4.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0958'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.43%  '
$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.89'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.477.07'
$ws.Range("E45").Value = '  +2.16%  '
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.99'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.40%  '
$ws.Range("E49").Value = '  -1.75%  '
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("D51").Value = '2.252.11'
$ws.Range("E51").Value = '  -1.07%  '
